# May 9th changes: prepend 15 new rows of sample data above the existing
# data set, dropping what used to be the last 5 rows so the sheet still
# ends up with 30 data rows (A1:C31) total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows 2-21 down by 15 rows (bottom-up so we never
# overwrite a source row before it has been copied). This lands the old
# data in rows 17-36; rows 32-36 (the old rows 17-21) are removed afterward.
for ($r = 21; $r -ge 2; $r--) {
    $destRow = $r + 15
    $ws.Cells.Item($destRow, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($destRow, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}

# Populate the newly freed rows 2-16 with the new data added on May 9th.
$newData = New-Object 'object[,]' 15,3
$newData[0,0] = -0.0774271711707115
$newData[0,1] = 0.1702786833047866
$newData[0,2] = -0.0888808965682983
$newData[1,0] = -0.047036625444889
$newData[1,1] = 0.0433714315295219
$newData[1,2] = -0.0064140851609408
$newData[2,0] = 0.0061086523346602
$newData[2,1] = -0.0146607663482427
$newData[2,2] = -0.0610865242779254
$newData[3,0] = 0.0401643887162208
$newData[3,1] = 0.0106901414692401
$newData[3,2] = 0.0001527163112768
$newData[4,0] = 0.0146607663482427
$newData[4,1] = 0.0265726372599601
$newData[4,2] = 0.0464257597923278
$newData[5,0] = -0.0455094613134861
$newData[5,1] = 0.0074830991216003
$newData[5,2] = 0.0467311926186084
$newData[6,0] = -0.0215329993516206
$newData[6,1] = 0.0041233403608202
$newData[6,2] = 0.0065668015740811
$newData[7,0] = 0.0096211275085806
$newData[7,1] = -0.0128281703218817
$newData[7,2] = 0.0253509078174829
$newData[8,0] = 0.024892758578062
$newData[8,1] = -0.0320704244077205
$newData[8,2] = 0.0131336031481623
$newData[9,0] = 0.0187841057777404
$newData[9,1] = -0.0146607663482427
$newData[9,2] = -0.0226020142436027
$newData[10,0] = -0.028557950630784
$newData[10,1] = -0.0305432621389627
$newData[10,2] = -0.0201585534960031
$newData[11,0] = -0.0192422550171613
$newData[11,1] = 0.0343611687421798
$newData[11,2] = -0.0004581489483825
$newData[12,0] = 0.009010262787342
$newData[12,1] = 0.0630718395113945
$newData[12,2] = -0.0142026171088218
$newData[13,0] = 0.0442877300083637
$newData[13,1] = 0.0134390350431203
$newData[13,2] = -0.001527163083665
$newData[14,0] = -0.0262672062963247
$newData[14,1] = -0.027030786499381
$newData[14,2] = -0.0058032199740409

$ws.Range("A2:C16").Value2 = $newData

# The last 5 rows of the original data (old rows 17-21, now shifted down to
# rows 32-36) are not part of the final data set, so remove them.
$ws.Range("A32:C36").EntireRow.Delete()
